# Gantt chart update: add "2nd revision" review loop (t13/t14) and a new
# prep task (t15), refresh several dates / % complete figures, and fold a
# new dependency into t8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 (t2 "ID writing"): now fully done; drop the Priority flag ----
$ws.Range("F3").Value = 100
$ws.Range("H3").Clear()

# ---- Row 4 (t3 "IP filing"): now fully done; drop the Priority flag ----
$ws.Range("F4").Value = 100
$ws.Range("H4").Clear()

# ---- Row 6 (t4 "JDA GLW revision to BOE"): duration updated ----
$ws.Range("E6").Value = 31

# ---- Row 7 (t10 "BOE revision to GLW"): duration + % complete added ----
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 100

# ---- Row 8 (t11 "GLW final review"): shorter duration, now depends on new t14 ----
$ws.Range("E8").Value = 5
$ws.Range("I8").Value = "t14"

# ---- Row 9 (t12 "BOE final review"): duration updated ----
$ws.Range("E9").Value = 10

# ---- Row 11 (t6 "BLU components shipping to SP"): new date, longer, partially done ----
$ws.Range("C11").Value = "May 10, 2020"
$ws.Range("C11").HorizontalAlignment = -4108  # xlCenter (fix stray right-align)
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 80

# ---- Row 12 (t8 "Prepare PGD, assemble BLU"): now also depends on new t15 ----
$ws.Range("B12").HorizontalAlignment = -4108  # xlCenter (fix stray right-align)
$ws.Range("I12").Value = "t6,t15"
$ws.Range("I12").HorizontalAlignment = -4152  # xlRight

# ---- New row 16: t13 "GLW 2nd revision", depends on t10 ----
$ws.Range("A16").Value = "t13"
$ws.Range("B16").Value = "GLW 2nd revision"
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 100
$ws.Range("G16").Value = "Ryan"
$ws.Range("I16").Value = "t10"
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("G16").HorizontalAlignment = -4108
$ws.Range("I16").HorizontalAlignment = -4108

# ---- New row 17: t14 "BOE 2nd revision", depends on t13 ----
$ws.Range("A17").Value = "t14"
$ws.Range("B17").Value = "BOE 2nd revision"
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100
$ws.Range("G17").Value = "Tom"
$ws.Range("I17").Value = "t13"
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("G17").HorizontalAlignment = -4108
$ws.Range("I17").HorizontalAlignment = -4108

# ---- New row 18: t15 "Get lab access and prework done" ----
$ws.Range("A18").Value = "t15"
$ws.Range("B18").Value = "Get lab access and prework done"
$ws.Range("C18").Value = "Apr 15, 2020"
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 100
$ws.Range("G18").Value = "Mi"
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("C18").HorizontalAlignment = -4108
$ws.Range("G18").HorizontalAlignment = -4108
